# Update "Sprint 3" backlog sheet - fill in the second Thursday/Friday
# (columns P and Q) actual-remaining-hours data for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 3")

$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 1

$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 1

$ws.Range("Q4").Value = 0

$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0

$excel.CalculateFullRebuild()
